$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at the very left; existing columns A-E shift to B-F
$ws.Columns.Item(1).Insert()

# Rename old "ID" header (now shifted into column B) to "MovieId", then set new column A header
$ws.Range("B1").Value = "MovieId"
$ws.Range("A1").Value = "id"

# Fill in the new GUID identifiers for each movie row
$ws.Range("A2").Value = "a5521c81-8a9e-4ac5-8031-52008277c4ec"
$ws.Range("A3").Value = "62ae2134-1ad3-4496-af69-e86318abb836"
$ws.Range("A4").Value = "b21ea323-8d40-42c7-91f3-213e5dbfef55"
$ws.Range("A5").Value = "ae53c082-b231-4dab-9647-900d066eeed8"
$ws.Range("A6").Value = "1a6b4927-ac51-4453-8d82-5a7227511e09"

# Auto-fit the new id column so its width matches the data (GUIDs)
$ws.Columns.Item(1).AutoFit()

# Match the new selection shown in the diff
[void]$ws.Range("A6").Select()
